# Insert a new weekly price-report row for Cereza "Early Burlat" at row 134
# (pushing the existing rows 134:225 down to 135:226).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(134).Insert()

$ws.Cells.Item(134, 1).Value = 10
$ws.Cells.Item(134, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(134, 3).Value = "La Araucanía"
$ws.Cells.Item(134, 4).Value = 44879
$ws.Cells.Item(134, 5).Value = 9
$ws.Cells.Item(134, 6).Value = "Fruta"
$ws.Cells.Item(134, 7).Value = 100103
$ws.Cells.Item(134, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(134, 9).Value = 100103001
$ws.Cells.Item(134, 10).Value = "Cereza"
$ws.Cells.Item(134, 11).Value = "Early Burlat"
$ws.Cells.Item(134, 12).Value = "Especial"
$ws.Cells.Item(134, 13).Value = 75
$ws.Cells.Item(134, 14).Value = 25000
$ws.Cells.Item(134, 15).Value = 25000
$ws.Cells.Item(134, 16).Value = 25000
$ws.Cells.Item(134, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(134, 18).Value = "Región Metropolitana"
$ws.Cells.Item(134, 19).Value = 2500
$ws.Cells.Item(134, 20).Value = 10
